$d = $word.ActiveDocument

# The paragraph ending in the "petascale@shodor.org" hyperlink is immediately
# followed by an otherwise-empty paragraph that only contains a page-break run.
# Style consistency: merge them into a single paragraph (delete the intervening
# paragraph mark), then move the "_GoBack" bookmark that currently sits a few
# paragraphs further down (right after the "top500.org" paragraph) up to sit
# right after the hyperlink, before the page-break run.

$mailParagraph = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -match "petascale@shodor\.org") {
        $mailParagraph = $p
    }
}

# Delete the paragraph mark at the end of this paragraph, merging it with the
# following (page-break-only) paragraph.
$pMark = $d.Range($mailParagraph.Range.End - 1, $mailParagraph.Range.End)
$pMark.Delete()

# Locate the end of the "petascale@shodor.org" text so we can drop the
# bookmark right after it (and before the page-break run that now follows it
# in the same paragraph).
$found = $d.Content
$found.Find.Execute("petascale@shodor.org", $false, $false, $false, $false, $false, `
                     $true, 1, $false, "", 0)
$bmRange = $d.Range($found.End, $found.End)

# Re-adding a bookmark named "_GoBack" removes any pre-existing "_GoBack"
# bookmark elsewhere in the document (Word only ever keeps a single instance),
# so this both relocates it and cleans up the old empty-paragraph occurrence.
$d.Bookmarks.Add("_GoBack", $bmRange)
